# Weekly price-sheet update: a new record (Fecha=45180, Volumen=50) is added
# at the top of the data table (row 369), pushing the existing rows 369-425
# down by one (to 370-426). The new row keeps the same Mercado/Región/
# Categoría/etc. values as the row it displaced, only the date (column D)
# and the volume (column J) are new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 369; this shifts rows 369:425 down to 370:426
# and grows the sheet's used range / dimension accordingly.
$ws.Rows.Item(369).Insert()

# The row that used to be 369 is now at 370 - clone its values into the
# freshly inserted row 369.
$srcRow = $ws.Range("A370:R370")
$newRow = $ws.Range("A369:R369")
$newRow.Value2 = $srcRow.Value2

# Match the date cell's number format (column D uses a date/time format).
$ws.Cells.Item(369, 4).NumberFormat = $ws.Cells.Item(370, 4).NumberFormat

# Overwrite the two cells that actually carry new data for this entry.
$ws.Cells.Item(369, 4).Value2 = 45180   # Fecha
$ws.Cells.Item(369, 10).Value2 = 50     # Volumen
